$d = $word.ActiveDocument

# The document contains four "<id>...</id>" markers (p065r_1..p065r_4),
# each originally split across three runs:
#   run1: "<id>"      (Courier New, color 7f6000, sz 18)
#   run2: "p065r_N"    (color 000000)
#   run3: "</id>"     (Courier New, color 7f6000, sz 18)
# They need to be merged into a single run (keeping run1's formatting)
# whose text reads "<id>p065r_N</id>". The sibling "fig_p065r_N" markers
# must stay untouched.

$ids = @("p065r_1", "p065r_2", "p065r_3", "p065r_4")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"

    # Locate the exact marker text in the document.
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($found) {
        $start = $rng.Start
        $end = $rng.End
        $firstRunLen = 4   # length of "<id>"
        $firstRunEnd = $start + $firstRunLen

        # Remove the text belonging to run2 ("p065r_N") and run3 ("</id>"),
        # leaving only run1 ("<id>") behind.
        $tail = $d.Range($firstRunEnd, $end)
        $tail.Text = ""

        # Re-append the removed content to run1 so it becomes one run
        # carrying run1's original character formatting.
        $head = $d.Range($start, $firstRunEnd)
        $head.InsertAfter($id + "</id>")
    }
}
